# Incluimos la versión pdf
#
# Fixes a handful of small text/typo issues across the deck:
#  - Slide 12 title: "Resumen I" -> "Resumen"
#  - Slide 12 body: "implentar" -> "implementar" (typo fix)
#  - Slide 6 body: hyperlink run "link" -> "Link" (capitalisation)
#  - Slide 8 body: hyperlink run "link" -> "Link" (capitalisation)

$p = $ppt.ActivePresentation

# --- Slide 12: "Resumen I" -> "Resumen" -------------------------------
$slide12 = $p.Slides.Item(12)
$title12 = $slide12.Shapes.Item(1)
$title12.TextFrame.TextRange.Text = "Resumen"

# --- Slide 12: fix typo "implentar" -> "implementar" -------------------
$body12 = $slide12.Shapes.Item(2)
$body12Text = $body12.TextFrame.TextRange
$lastPara = $body12Text.Paragraphs($body12Text.Paragraphs().Count)
$lastPara.Runs(1).Text = "Son bastante sencillas de implementar y se aportan varios ejemplos funcionales en las internas."

# --- Slide 6: capitalize hyperlink text "link" -> "Link" ---------------
$slide6 = $p.Slides.Item(6)
$body6 = $slide6.Shapes.Item(2)
$body6Text = $body6.TextFrame.TextRange
$body6LastPara = $body6Text.Paragraphs($body6Text.Paragraphs().Count)
$body6LastPara.Runs($body6LastPara.Runs().Count).Text = "Link"

# --- Slide 8: capitalize hyperlink text "link" -> "Link" ---------------
$slide8 = $p.Slides.Item(8)
$body8 = $slide8.Shapes.Item(2)
$body8Para1 = $body8.TextFrame.TextRange.Paragraphs(1)
$body8Para1.Runs(2).Text = "Link"
